$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.466.79'
$ws.Range("E2").Value = '  +4.97%  '
$ws.Range("D3").Value = '2.491.24'
$ws.Range("E3").Value = '  +2.65%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '322.58'
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("D6").Value = '105.54'
$ws.Range("E6").Value = '  +3.30%  '
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("D10").Value = '37.86'
$ws.Range("E10").Value = '  +6.59%  '
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '18.29'
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").Value = '7.15'
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '2.885.51'
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").Value = '2.499.77'
$ws.Range("E16").Value = '  +2.82%  '
$ws.Range("D17").Value = '0.843'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = '47.341.80'
$ws.Range("E18").Value = '  +4.90%  '
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("D20").Value = '6.55'
$ws.Range("E20").Value = '  +3.09%  '
$ws.Range("D21").Value = '0.0₃0935'
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("D22").Value = '70.68'
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").Value = '250.49'
$ws.Range("E23").Value = '  +2.81%  '
$ws.Range("E24").Value = '  +5.56%  '
$ws.Range("D25").Value = '2.56'
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").Value = '26.18'
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("E28").Value = '  +4.54%  '
$ws.Range("E29").Value = '  +6.68%  '
$ws.Range("D30").Value = '34.99'
$ws.Range("E30").Value = '  +6.26%  '
$ws.Range("E31").Value = '  +7.01%  '
$ws.Range("D32").Value = '49.48'
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("E33").Value = '  -2.03%  '
$ws.Range("D34").Value = '5.33'
$ws.Range("E34").Value = '  +2.63%  '
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").Value = '1.94'
$ws.Range("E37").Value = '  +3.57%  '
$ws.Range("D38").Value = '4.64'
$ws.Range("E38").Value = '  +4.19%  '
$ws.Range("E39").Value = '  +4.74%  '
$ws.Range("D40").Value = '2.25'
$ws.Range("E40").Value = '  +1.80%  '
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("E42").Value = '  -3.32%  '
$ws.Range("D43").Value = '20.90'
$ws.Range("E43").Value = '  +1.33%  '
$ws.Range("D44").Value = '0.0297'
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("D45").Value = '1.961.65'
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("D46").Value = '2.98'
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").Value = '9.21'
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").Value = '1.79'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").Value = '5.31'
$ws.Range("E50").Value = '  +12.51%  '
$ws.Range("D51").Value = '79.16'
$ws.Range("E51").Value = '  +3.62%  '
